$d = $word.ActiveDocument

# Merge "github.com/" + "somguynamedseb" into one run (drop the spell-check
# proofErr wrapper in the process) -- textually this is a no-op replace.
$d.Content.Find.Execute("github.com/somguynamedseb", $true, $false, $false, $false, $false,
                         $true, 1, $false, "github.com/somguynamedseb", 2)

# Replace the github-pages URL with the personal domain, prefixed with the
# same amount of leading whitespace that used to sit between the two runs.
$d.Content.Find.Execute("somguynamedseb.github.io", $true, $false, $false, $false, $false,
                         $true, 1, $false, "      sebastianbaldini.com", 2)

# Merge "linkedin.com/in/" + "sebastian-baldini" into one run (drop the
# spell-check proofErr wrapper in the process).
$d.Content.Find.Execute("linkedin.com/in/sebastian-baldini", $true, $false, $false, $false, $false,
                         $true, 1, $false, "linkedin.com/in/sebastian-baldini", 2)
